$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.962.70'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '3.385.48'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.48'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.15'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("D12").Value = '3.964.82'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.78'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").Value = '3.376.39'
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").Value = '61.074.86'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.08'
$ws.Range("E18").Value = '  -3.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.66'
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.98'
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.48'
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("E22").Value = '  +2.74%  '
$ws.Range("E23").Value = '  -2.45%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  -4.72%  '
$ws.Range("D26").Value = '3.524.81'
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.00'
$ws.Range("E30").Value = '  -2.10%  '
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.39'
$ws.Range("E33").Value = '  -4.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.27'
$ws.Range("E34").Value = '  -2.30%  '
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '166.85'
$ws.Range("D37").Value = '3.416.94'
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("E39").Value = '  -4.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0767'
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.62'
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("E44").Value = '  -2.34%  '
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.13'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '2.457.59'
$ws.Range("E47").Value = '  -4.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.02'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("E50").Value = '  +7.23%  '
$ws.Range("E51").Value = '  +1.23%  '
